# Weekly update: insert two new daily price records (rows 106-107) for
# "Arveja Verde" at Mercado Mayorista Lo Valledor de Santiago. All the
# existing rows from 106 downward shift down by two (to 108-171), which
# Rows.Insert() handles for us automatically (including extending the
# sheet's used-range dimension from A1:R169 to A1:R171).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 106:169 down to 108:171, leaving two fresh blank
# rows at 106:107 for the new records.
$ws.Rows("106:107").Insert()

# New row 106
$ws.Range("A106").Value = 6
$ws.Range("B106").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C106").Value = 'Metropolitana'
$ws.Range("D106").Value = 44518
$ws.Range("E106").Value = 13
$ws.Range("F106").Value = 100112022
$ws.Range("G106").Value = 'Arveja Verde'
$ws.Range("H106").Value = 'Sin especificar'
$ws.Range("I106").Value = 'Primera'
$ws.Range("J106").Value = 350
$ws.Range("K106").Value = 15000
$ws.Range("L106").Value = 17000
$ws.Range("M106").Value = 15857
$ws.Range("N106").Value = '$/saco 25 kilos'
$ws.Range("O106").Value = 'Carahue'
$ws.Range("P106").Value = 634
$ws.Range("Q106").Value = 25
$ws.Range("R106").Value = 'Hortaliza'

# New row 107
$ws.Range("A107").Value = 6
$ws.Range("B107").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C107").Value = 'Metropolitana'
$ws.Range("D107").Value = 44518
$ws.Range("E107").Value = 13
$ws.Range("F107").Value = 100112022
$ws.Range("G107").Value = 'Arveja Verde'
$ws.Range("H107").Value = 'Sin especificar'
$ws.Range("I107").Value = 'Primera'
$ws.Range("J107").Value = 600
$ws.Range("K107").Value = 14000
$ws.Range("L107").Value = 15000
$ws.Range("M107").Value = 14333
$ws.Range("N107").Value = '$/saco 25 kilos'
$ws.Range("O107").Value = 'Región del Maule'
$ws.Range("P107").Value = 573
$ws.Range("Q107").Value = 25
$ws.Range("R107").Value = 'Hortaliza'
